$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (e.g. "1.000", "238.29") need the cell pre-formatted as Text so the literal
# string is preserved exactly, matching the source data (inline strings).
$textCells = @(
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D13",
    "D14",
    "D15",
    "D17",
    "D18",
    "D19",
    "D20",
    "D21",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D38",
    "D39",
    "D40",
    "D42",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.298.45"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.859.71"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "0.7038"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value = "238.29"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.07888"
$ws.Range("E8").Value = "  +2.28%  "
$ws.Range("D9").Value = "0.3037"
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("D10").Value = "24.40"
$ws.Range("D11").Value = "0.08183"
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "1.872.88"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").Value = "0.7200"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").Value = "5.214"
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").Value = "89.67"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "29.317.13"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "5.808"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").Value = "0.000007822"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").Value = "13.25"
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("D20").Value = "238.05"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.114.80"
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "7.560"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("D25").Value = "162.13"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("D26").Value = "8.902"
$ws.Range("E26").Value = "  -1.41%  "
$ws.Range("D27").Value = "0.1421"
$ws.Range("E27").Value = "  -4.58%  "
$ws.Range("D28").Value = "18.11"
$ws.Range("D29").Value = "1.920"
$ws.Range("E29").Value = "  -4.05%  "
$ws.Range("D30").Value = "1.390"
$ws.Range("E30").Value = "  -1.77%  "
$ws.Range("D31").Value = "1.478"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("D32").Value = "4.323"
$ws.Range("E32").Value = "  -3.10%  "
$ws.Range("D33").Value = "4.052"
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("D34").Value = "0.05187"
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("D35").Value = "1.175"
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("D36").Value = "0.7151"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "2.677"
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("D39").Value = "0.01850"
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("D40").Value = "2.687"
$ws.Range("E40").Value = "  -1.54%  "
$ws.Range("D41").Value = "1.151.04"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("D42").Value = "0.9228"
$ws.Range("E42").Value = "  -1.47%  "
$ws.Range("E43").Value = "  +1.29%  "
$ws.Range("D44").Value = "0.4257"
$ws.Range("E44").Value = "  -1.07%  "
$ws.Range("D45").Value = "70.79"
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "102.00"
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("D48").Value = "0.5318"
$ws.Range("E48").Value = "  -3.18%  "
$ws.Range("D49").Value = "1.754"
$ws.Range("E49").Value = "  -2.58%  "
$ws.Range("D50").Value = "9.179"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "7.005"
$ws.Range("E51").Value = "  -0.02%  "
